# Add a new bulleted list item after "Viết manual instructions." with the
# same list paragraph style/numbering, per the author's note about
# BadgeChecker tokens.

$d = $word.ActiveDocument

# Locate the last paragraph in the document (the "Viết manual instructions."
# bullet) and append a brand-new paragraph right after it, inheriting the
# same paragraph formatting (ListParagraph style + list numbering).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphAfter()

# The freshly inserted paragraph is now the new last paragraph; set its text.
$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.Text = "Sau một số lần làm việc, thì nó được tặng token, token này có thể để mua đồ thiết kế huy hiệu."
